# Rename several header labels across all worksheets of the workbook.
# Sheets 1-8  ("LLL_*")            -> short header row A1:Q1
# Sheets 9-32 ("LL_*","LLG_*","LG_*") -> long header row A1:AQ1
#
# Renames applied (old -> new):
#   pf_ikss_from_deg       -> pf_ikss_from_degree      (col L, sheets 1-8)
#   pf_ikss_to_deg         -> pf_ikss_to_degree        (col M, sheets 1-8)
#   pf_va_from_deg         -> pf_va_from_degree        (col P, sheets 1-8)
#   pf_va_to_deg           -> pf_va_to_degree          (col Q, sheets 1-8)
#
#   pf_q_a_from_mw         -> pf_q_a_from_mvar         (col T,  sheets 9-32)
#   pf_q_b_from_mw         -> pf_q_b_from_mvar         (col U,  sheets 9-32)
#   pf_q_c_from_mw         -> pf_q_c_from_mvar         (col V,  sheets 9-32)
#   pf_q_a_to_mw           -> pf_q_a_to_mvar           (col W,  sheets 9-32)
#   pf_q_b_to_mw           -> pf_q_b_to_mvar           (col X,  sheets 9-32)
#   pf_q_c_to_mw           -> pf_q_c_to_mvar           (col Y,  sheets 9-32)
#   pf_ikss_a_from_deg     -> pf_ikss_a_from_degree    (col Z,  sheets 9-32)
#   pf_ikss_b_from_deg     -> pf_ikss_b_from_degree    (col AA, sheets 9-32)
#   pf_ikss_c_from_deg     -> pf_ikss_c_from_degree    (col AB, sheets 9-32)
#   pf_ikss_a_to_deg       -> pf_ikss_a_to_degree      (col AC, sheets 9-32)
#   pf_ikss_b_to_deg       -> pf_ikss_b_to_degree      (col AD, sheets 9-32)
#   pf_ikss_c_to_deg       -> pf_ikss_c_to_degree      (col AE, sheets 9-32)
#   pf_vm_b_from_bus_pu    -> pf_vm_b_from_pu          (col AG, sheets 9-32)
#   pf_vm_c_from_bus_pu    -> pf_vm_c_from_pu          (col AH, sheets 9-32)
#   pf_vm_a_to_bus_pu      -> pf_vm_a_to_pu            (col AI, sheets 9-32)
#   pf_vm_b_to_bus_pu      -> pf_vm_b_to_pu            (col AJ, sheets 9-32)
#   pf_vm_c_to_bus_pu      -> pf_vm_c_to_pu            (col AK, sheets 9-32)
#   pf_va_a_from_bus_deg   -> pf_va_a_from_degree      (col AL, sheets 9-32)
#   pf_va_b_from_bus_deg   -> pf_va_b_from_degree      (col AM, sheets 9-32)
#   pf_va_c_from_bus_deg   -> pf_va_c_from_degree      (col AN, sheets 9-32)
#   pf_va_a_to_bus_deg     -> pf_va_a_to_degree        (col AO, sheets 9-32)
#   pf_va_b_to_bus_deg     -> pf_va_b_to_degree        (col AP, sheets 9-32)
#   pf_va_c_to_bus_deg     -> pf_va_c_to_degree        (col AQ, sheets 9-32)

$wb = $excel.ActiveWorkbook

# Mapping applied to every worksheet that has a short header row (A1:Q1)
$shortHeaderMap = @{
    "pf_ikss_from_deg" = "pf_ikss_from_degree"
    "pf_ikss_to_deg"   = "pf_ikss_to_degree"
    "pf_va_from_deg"   = "pf_va_from_degree"
    "pf_va_to_deg"     = "pf_va_to_degree"
}

# Mapping applied to every worksheet that has a long header row (A1:AQ1)
$longHeaderMap = @{
    "pf_q_a_from_mw"       = "pf_q_a_from_mvar"
    "pf_q_b_from_mw"       = "pf_q_b_from_mvar"
    "pf_q_c_from_mw"       = "pf_q_c_from_mvar"
    "pf_q_a_to_mw"         = "pf_q_a_to_mvar"
    "pf_q_b_to_mw"         = "pf_q_b_to_mvar"
    "pf_q_c_to_mw"         = "pf_q_c_to_mvar"
    "pf_ikss_a_from_deg"   = "pf_ikss_a_from_degree"
    "pf_ikss_b_from_deg"   = "pf_ikss_b_from_degree"
    "pf_ikss_c_from_deg"   = "pf_ikss_c_from_degree"
    "pf_ikss_a_to_deg"     = "pf_ikss_a_to_degree"
    "pf_ikss_b_to_deg"     = "pf_ikss_b_to_degree"
    "pf_ikss_c_to_deg"     = "pf_ikss_c_to_degree"
    "pf_vm_b_from_bus_pu"  = "pf_vm_b_from_pu"
    "pf_vm_c_from_bus_pu"  = "pf_vm_c_from_pu"
    "pf_vm_a_to_bus_pu"    = "pf_vm_a_to_pu"
    "pf_vm_b_to_bus_pu"    = "pf_vm_b_to_pu"
    "pf_vm_c_to_bus_pu"    = "pf_vm_c_to_pu"
    "pf_va_a_from_bus_deg" = "pf_va_a_from_degree"
    "pf_va_b_from_bus_deg" = "pf_va_b_from_degree"
    "pf_va_c_from_bus_deg" = "pf_va_c_from_degree"
    "pf_va_a_to_bus_deg"   = "pf_va_a_to_degree"
    "pf_va_b_to_bus_deg"   = "pf_va_b_to_degree"
    "pf_va_c_to_bus_deg"   = "pf_va_c_to_degree"
}

foreach ($ws in $wb.Worksheets) {
    $headerRow = $ws.Range("A1:AQ1")
    foreach ($cell in $headerRow.Cells) {
        $val = $cell.Value2
        if ($null -ne $val) {
            if ($shortHeaderMap.ContainsKey($val)) {
                $cell.Value2 = $shortHeaderMap[$val]
            } elseif ($longHeaderMap.ContainsKey($val)) {
                $cell.Value2 = $longHeaderMap[$val]
            }
        }
    }
}
